# Update odds values in row 2 of the active worksheet to match the
# latest Betfair Back/Lay snapshot for 2025-12-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1.83
    "H2"  = 5
    "I2"  = 6
    "J2"  = 3.75
    "K2"  = 4.1
    "L2"  = 1.4
    "N2"  = 3.75
    "O2"  = 1.32
    "P2"  = 1.95
    "Q2"  = 1.98
    "R2"  = 1.36
    "S2"  = 3.55
    "T2"  = 1.92
    "U2"  = 1.92
    "W2"  = 2.2
    "X2"  = 14
    "Y2"  = 19
    "Z2"  = 42
    "AA2" = 150
    "AB2" = 8.6
    "AC2" = 9.2
    "AE2" = 80
    "AF2" = 10.5
    "AG2" = 10
    "AH2" = 22
    "AI2" = 80
    "AJ2" = 22
    "AK2" = 19
    "AO2" = 100
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
